$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, pushing the existing rows 139-189
# (and the sheet dimension) down by one, to 140-190.
$ws.Rows("139:139").Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(139, 1).Value  = 5
$ws.Cells.Item(139, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(139, 3).Value  = "Maule"
$ws.Cells.Item(139, 4).Value  = 44468
$ws.Cells.Item(139, 5).Value  = 7
$ws.Cells.Item(139, 6).Value  = 100114013
$ws.Cells.Item(139, 7).Value  = "Zanahoria"
$ws.Cells.Item(139, 8).Value  = "Sin especificar"
$ws.Cells.Item(139, 9).Value  = "Primera"
$ws.Cells.Item(139, 10).Value = 400
$ws.Cells.Item(139, 11).Value = 6000
$ws.Cells.Item(139, 12).Value = 6000
$ws.Cells.Item(139, 13).Value = 6000
$ws.Cells.Item(139, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(139, 15).Value = "Región de Ñuble"
$ws.Cells.Item(139, 16).Value = 300
$ws.Cells.Item(139, 17).Value = 20
$ws.Cells.Item(139, 18).Value = "Hortaliza"
